$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: date corrected from 45336 to 45335
$ws.Range("B21").Value = 45335

# Row 22: fill in a new time-log entry (Aris, 2023-... date 45335,
# start-time formula 14+32/60, duration formula recalculates automatically)
$ws.Range("A22").Value = "Aris"

# copy date formatting/border from B21 onto B22 before writing the value
$ws.Range("B21").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("B22").Value = 45335

$ws.Range("C22").Formula = "=14+32/60"

# move the active selection to C23 (result of the button-driven move)
$ws.Range("C23").Select()
